$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------------
# 1. Update the "time_taken" (column F) timestamps on the existing "data"
#    sheet (rows 2-49) to the refreshed panel_query_time run.
# ---------------------------------------------------------------------------
$newTimestamps = @(
  "2021-10-05 14:21:42.825751",
  "2021-10-05 14:21:42.825759",
  "2021-10-05 14:21:42.825762",
  "2021-10-05 14:21:42.825765",
  "2021-10-05 14:21:42.825768",
  "2021-10-05 14:21:42.825771",
  "2021-10-05 14:21:42.825774",
  "2021-10-05 14:21:42.825776",
  "2021-10-05 14:21:42.825779",
  "2021-10-05 14:21:42.825782",
  "2021-10-05 14:21:42.825784",
  "2021-10-05 14:21:42.825787",
  "2021-10-05 14:21:42.825790",
  "2021-10-05 14:21:42.825792",
  "2021-10-05 14:21:42.825795",
  "2021-10-05 14:21:42.825797",
  "2021-10-05 14:21:42.825800",
  "2021-10-05 14:21:42.825803",
  "2021-10-05 14:21:42.825805",
  "2021-10-05 14:21:42.825808",
  "2021-10-05 14:21:42.825810",
  "2021-10-05 14:21:42.825813",
  "2021-10-05 14:21:42.825815",
  "2021-10-05 14:21:42.825818",
  "2021-10-05 14:21:42.825821",
  "2021-10-05 14:21:42.825823",
  "2021-10-05 14:21:42.825826",
  "2021-10-05 14:21:42.825828",
  "2021-10-05 14:21:42.825831",
  "2021-10-05 14:21:42.825833",
  "2021-10-05 14:21:42.825836",
  "2021-10-05 14:21:42.825838",
  "2021-10-05 14:21:42.825841",
  "2021-10-05 14:21:42.825844",
  "2021-10-05 14:21:42.825847",
  "2021-10-05 14:21:42.825849",
  "2021-10-05 14:21:42.825852",
  "2021-10-05 14:21:42.825854",
  "2021-10-05 14:21:42.825857",
  "2021-10-05 14:21:42.825859",
  "2021-10-05 14:21:42.825862",
  "2021-10-05 14:21:42.825865",
  "2021-10-05 14:21:42.825867",
  "2021-10-05 14:21:42.825870",
  "2021-10-05 14:21:42.825872",
  "2021-10-05 14:21:42.825875",
  "2021-10-05 14:21:42.825878",
  "2021-10-05 14:21:42.825880"
)

for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
  $row = $i + 2
  $dataSheet.Cells.Item($row, 6).Value = $newTimestamps[$i]
}

# ---------------------------------------------------------------------------
# 2. Add the new "metadata" sheet, placed right after "data".
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$ws.Name = "metadata"

# Header row (bold / bordered / centered - matches the "data" sheet header
# style), copied from the existing header so it reuses the same style id.
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

$dataSheet.Range("B1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)

# Data row.
$ws.Range("A2").Value = 0
$dataSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("B2").Value = "Nephrocalcinosis or nephrolithiasis"
$ws.Range("C2").Value = 149

# data_version needs to stay text ("2.24"), not be coerced to a number.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2.24"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "2021-08-19T11:33:56.974648Z"
$ws.Range("F2").Value = "2021-10-05 14:21:42.822035"
$ws.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/149/?format=json"
